$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Eritrea" and "Islas Turcas y Caicos" in the shared-string order ---
# Row 180 currently shows "Eritrea", row 181 currently shows "Islas Turcas y Caicos".
# After the edit, row 180 should show "Islas Turcas y Caicos" and row 181 "Eritrea".
$ws.Range("A180").Value = "Islas Turcas y Caicos"
$ws.Range("A181").Value = "Eritrea"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 05:45"

# --- Kazajistan (row 29) ---
$ws.Range("B29").Value = 103571
$ws.Range("C29").Value = 271
$ws.Range("D29").Value = 86286
$ws.Range("E29").Value = 15870

# --- Belgica (row 40) ---
$ws.Range("B40").Value = 78897
$ws.Range("C40").Value = 363
$ws.Range("D40").Value = 18048
$ws.Range("E40").Value = 50890
$ws.Range("G40").Value = 15
$ws.Range("H40").Value = 9959

# --- San Martin (Parte Holandesa) (row 179) ---
$ws.Range("B179").Value = 333
$ws.Range("C179").Value = 7
$ws.Range("E179").Value = 209

# --- Islas Turcas y Caicos (now row 180, after the swap above) ---
$ws.Range("B180").Value = 315
$ws.Range("C180").Value = 17
$ws.Range("D180").Value = 56
$ws.Range("E180").Value = 257
$ws.Range("H180").Value = 2

# --- Eritrea (now row 181, after the swap above) ---
$ws.Range("B181").Value = 304
$ws.Range("D181").Value = 261
$ws.Range("E181").Value = 43
$ws.Range("H181").Value = 0

# --- Camboya (row 183) ---
$ws.Range("D183").Value = 251
$ws.Range("E183").Value = 22
